$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I1").Value = "price_after_discount"
$ws.Range("I1").Select()
